$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before current column M (13) so the "Model_Base" header
# and its data shift from M to N, making room for the new "NumOptions" column.
$ws.Columns.Item(13).Insert()

# New header for the inserted column M
$ws.Range("M2").Value = "NumOptions"

# New data values for the "NumOptions" column (rows 3-7)
$ws.Range("M3").Value = 4
$ws.Range("M4").Value = 4
$ws.Range("M5").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("M7").Value = 0

# Update selection to match the authored workbook state
$ws.Range("M8").Select()
